$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) cells to remain Text so dotted/zero-padded
# numeric-looking values are not reinterpreted as numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "31.117.46"
$ws.Range("E2").Value = "  +1.75%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.960.41"
$ws.Range("E3").Value = "  +2.08%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.05"
$ws.Range("E5").Value = "  +0.67%  "

$ws.Range("E6").Value = "  +0.14%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4908"
$ws.Range("E7").Value = "  +1.68%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2978"
$ws.Range("E8").Value = "  +2.71%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06836"
$ws.Range("E9").Value = "  +0.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.10"
$ws.Range("E10").Value = "  -1.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "106.65"
$ws.Range("E11").Value = "  -4.95%  "

$ws.Range("E12").Value = "  +2.46%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.934.03"
$ws.Range("E13").Value = "  +0.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.424"
$ws.Range("E14").Value = "  -1.24%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7120"
$ws.Range("E15").Value = "  +5.72%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "285.87"
$ws.Range("E16").Value = "  -3.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "31.134.19"
$ws.Range("E17").Value = "  +1.80%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007796"
$ws.Range("E18").Value = "  +1.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.26"
$ws.Range("E19").Value = "  +1.69%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.17%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.564"
$ws.Range("E21").Value = "  +0.79%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.188.08"
$ws.Range("E22").Value = "  +1.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.600"
$ws.Range("E24").Value = "  +2.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.987"
$ws.Range("E25").Value = "  +5.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.85"
$ws.Range("E26").Value = "  +1.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.07"
$ws.Range("E27").Value = "  -1.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.202"
$ws.Range("E28").Value = "  +5.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1061"
$ws.Range("E29").Value = "  -0.32%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.443"
$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.790"
$ws.Range("E31").Value = "  +18.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.536"
$ws.Range("E32").Value = "  +9.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05022"
$ws.Range("E33").Value = "  +0.74%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7708"
$ws.Range("E34").Value = "  +4.92%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.168"
$ws.Range("E35").Value = "  +2.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02056"
$ws.Range("E36").Value = "  +1.33%  "

$ws.Range("E37").Value = "  +0.80%  "

$ws.Range("E38").Value = "  +1.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.139"
$ws.Range("E39").Value = "  +5.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.431"
$ws.Range("E40").Value = "  +9.68%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8858"
$ws.Range("E41").Value = "  +1.88%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4478"
$ws.Range("E42").Value = "  +0.89%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "109.57"
$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "73.41"
$ws.Range("E44").Value = "  +5.66%  "

$ws.Range("E45").Value = "  +0.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.025.27"
$ws.Range("E46").Value = "  +21.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.493"
$ws.Range("E47").Value = "  +3.24%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.423"
$ws.Range("E48").Value = "  +2.30%  "

$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1272"
$ws.Range("E49").Value = "  +3.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.03"
$ws.Range("E50").Value = "  +3.37%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2586"
$ws.Range("E51").Value = "  +3.02%  "

